# Update data parsing logic: append one new sensor reading (row 51) to each
# of the 4 worksheets, mirroring the existing row layout (time, total
# length, ID, actual length, checksum, and their *_DEC counterparts).

$wb = $excel.ActiveWorkbook

function Add-SensorRow($Sheet, $Row, $Time, $TotalLen, $Id, $ActualLen, $Checksum, $TotalLenDec, $IdDec, $ActualLenDec, $ChecksumDec) {
    $Sheet.Cells.Item($Row, 1).Value = $Time
    $Sheet.Cells.Item($Row, 2).Value = $TotalLen
    $Sheet.Cells.Item($Row, 3).Value = $Id
    $Sheet.Cells.Item($Row, 4).Value = $ActualLen
    $Sheet.Cells.Item($Row, 5).Value = $Checksum
    $Sheet.Cells.Item($Row, 6).Value = $TotalLenDec

    # ID_DEC is a 24-digit integer literal; force text formatting first so
    # it round-trips exactly instead of being coerced into a lossy double,
    # then clear the formatting again so the cell keeps the sheet's default
    # (unstyled) look.
    $Sheet.Cells.Item($Row, 7).NumberFormat = "@"
    $Sheet.Cells.Item($Row, 7).Value = $IdDec
    $Sheet.Cells.Item($Row, 7).ClearFormats()

    $Sheet.Cells.Item($Row, 8).Value = $ActualLenDec
    $Sheet.Cells.Item($Row, 9).Value = $ChecksumDec
}

# ROW35-FE-LIFTER
$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
Add-SensorRow $ws1 51 "2025-03-06 10:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

# ROW35-MID-LIFTER
$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
Add-SensorRow $ws2 51 "2025-03-06 10:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

# ROW02-FE-LIFTER
$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
Add-SensorRow $ws3 51 "2025-03-06 10:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

# ROW02-MID-LIFTER
$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
Add-SensorRow $ws4 51 "2025-03-06 10:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3
